$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cells (row 5) - bold style matching K5/L5
$ws.Range("M5").Value = "AlarmLoadingDetail"
$ws.Range("N5").Value = "StandbyLoadingDetail"

# Data cells rows 6-8 matching the styling of H6:H8 (left aligned)
$ws.Range("M6").Value = "Battery Alarm (A)"
$ws.Range("N6").Value = "Battery Standby (A)"

$ws.Range("M7").Value = "Battery Alarm (A)"
$ws.Range("N7").Value = "Battery Standby (A)"

$ws.Range("M8").Value = "Battery Alarm (A)"
$ws.Range("N8").Value = "Battery Standby (A)"

# Copy styles from existing reference cells to new columns
$ws.Range("A7").Copy()
$ws.Range("M5:N5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H6").Copy()
$ws.Range("M6:N8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update selection/view to match diff (scroll so column I is leftmost, select M7:N8)
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M7:N8").Select()
